$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 239.92857
$ws.Range("I2").Value = 235.3077
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 235.3077
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -122.3077
$ws.Range("N2").Value = -526
# Row 17
$ws.Range("H17").Value = 2506062
$ws.Range("J17").Value = 2506062
$ws.Range("L17").Value = 7518186
$ws.Range("N17").Value = -7518522
# Row 112
$ws.Range("H112").Value = 1110.3529
$ws.Range("J112").Value = 1107.0416
$ws.Range("L112").Value = 3321.1248
$ws.Range("N112").Value = -5537.1248
# Row 132
$ws.Range("H132").Value = 3061.8
$ws.Range("I132").Value = 2933.8635
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8801.5905
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6271.5905
$ws.Range("N132").Value = -17060
# Row 135
$ws.Range("H135").Value = 812.04346
$ws.Range("I135").Value = 679.875
$ws.Range("J135").Value = 1114.1428
$ws.Range("K135").Value = 6118.875
$ws.Range("L135").Value = 10027.2852
$ws.Range("M135").Value = -3583.875
$ws.Range("N135").Value = -15097.2852
# Row 137
$ws.Range("H137").Value = 2931.5625
$ws.Range("I137").Value = 3300.5
$ws.Range("J137").Value = 2808.5833
$ws.Range("K137").Value = 9901.5
$ws.Range("L137").Value = 8425.749899999999
$ws.Range("M137").Value = -7351.5
$ws.Range("N137").Value = -13525.7499
# Row 138
$ws.Range("H138").Value = 2946.058
$ws.Range("I138").Value = 2201.9167
$ws.Range("J138").Value = 3234.1128
$ws.Range("K138").Value = 6605.750100000001
$ws.Range("L138").Value = 9702.338400000001
$ws.Range("M138").Value = -1465.750100000001
$ws.Range("N138").Value = -19982.3384

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4950
$ws.Range("I61").Value = 6900
$ws.Range("K61").Value = 6900
$ws.Range("M61").Value = -6688
# Row 74
$ws.Range("H74").Value = 17253.25
$ws.Range("J74").Value = 4002
$ws.Range("L74").Value = 4002
$ws.Range("N74").Value = -5750
# Row 77
$ws.Range("H77").Value = 17253.25
$ws.Range("J77").Value = 4002
$ws.Range("L77").Value = 20010
$ws.Range("N77").Value = -28746
# Row 132
$ws.Range("H132").Value = 2569.476
$ws.Range("I132").Value = 1680
$ws.Range("K132").Value = 5040
$ws.Range("M132").Value = -2510
# Row 136
$ws.Range("H136").Value = 4950
$ws.Range("I136").Value = 6900
$ws.Range("K136").Value = 20700
$ws.Range("M136").Value = -18150

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1664.303
$ws.Range("I134").Value = 1396.92
$ws.Range("J134").Value = 2499.875
$ws.Range("K134").Value = 4190.76
$ws.Range("L134").Value = 7499.625
$ws.Range("M134").Value = -1655.76
$ws.Range("N134").Value = -12569.625
# Row 135
$ws.Range("H135").Value = 40780
$ws.Range("J135").Value = 40780
$ws.Range("L135").Value = 40780
$ws.Range("N135").Value = -50920
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
# Row 140
$ws.Range("H140").Value = 44917.65
$ws.Range("J140").Value = 44917.65
$ws.Range("L140").Value = 44917.65
$ws.Range("N140").Value = -55277.65

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1838.75
$ws.Range("I31").Value = 960.8
$ws.Range("J31").Value = 3887.3
$ws.Range("K31").Value = 960.8
$ws.Range("L31").Value = 3887.3
$ws.Range("M31").Value = -665.8
$ws.Range("N31").Value = -4477.3
# Row 34
$ws.Range("H34").Value = 1838.75
$ws.Range("I34").Value = 960.8
$ws.Range("J34").Value = 3887.3
$ws.Range("K34").Value = 960.8
$ws.Range("L34").Value = 3887.3
$ws.Range("M34").Value = -758.8
$ws.Range("N34").Value = -4291.3
# Row 58
$ws.Range("H58").Value = 3856.2307
$ws.Range("I58").Value = 3683.6667
$ws.Range("J58").Value = 4004.1428
$ws.Range("K58").Value = 3683.6667
$ws.Range("L58").Value = 4004.1428
$ws.Range("M58").Value = -3480.6667
$ws.Range("N58").Value = -4410.1428
# Row 134
$ws.Range("H134").Value = 3033.2632
$ws.Range("I134").Value = 3314.5334
$ws.Range("J134").Value = 1978.5
$ws.Range("K134").Value = 9943.600199999999
$ws.Range("L134").Value = 5935.5
$ws.Range("M134").Value = -7408.600199999999
$ws.Range("N134").Value = -11005.5
# Row 136
$ws.Range("H136").Value = 3856.2307
$ws.Range("I136").Value = 3683.6667
$ws.Range("J136").Value = 4004.1428
$ws.Range("K136").Value = 11051.0001
$ws.Range("L136").Value = 12012.4284
$ws.Range("M136").Value = -8501.000100000001
$ws.Range("N136").Value = -17112.4284
# Row 137
$ws.Range("H137").Value = 37780
$ws.Range("J137").Value = 37780
$ws.Range("L137").Value = 37780
$ws.Range("N137").Value = -47980
# Row 138
$ws.Range("H138").Value = 36572.145
$ws.Range("J138").Value = 36572.145
$ws.Range("L138").Value = 36572.145
$ws.Range("N138").Value = -46852.145
# Row 139
$ws.Range("H139").Value = 44500
$ws.Range("J139").Value = 44500
$ws.Range("L139").Value = 44500
$ws.Range("N139").Value = -54780
# Row 140
$ws.Range("H140").Value = 57294.855
$ws.Range("J140").Value = 57294.855
$ws.Range("L140").Value = 57294.855
$ws.Range("N140").Value = -67654.85500000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 28.666666
$ws.Range("I2").Value = 66.333336
$ws.Range("J2").Value = 9.833333
$ws.Range("K2").Value = 398.000016
$ws.Range("L2").Value = 58.999998
$ws.Range("M2").Value = -285.000016
$ws.Range("N2").Value = -284.999998
# Row 23
$ws.Range("H23").Value = 96.25
$ws.Range("I23").Value = 55.333332
$ws.Range("J23").Value = 105.69231
$ws.Range("K23").Value = 165.999996
$ws.Range("L23").Value = 317.07693
$ws.Range("M23").Value = 69.00000399999999
$ws.Range("N23").Value = -787.0769299999999
# Row 33
$ws.Range("H33").Value = 402.75
$ws.Range("I33").Value = 305.5
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 1833
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = -1550
$ws.Range("N33").Value = -3566
# Row 38
$ws.Range("H38").Value = 33505.8
$ws.Range("I38").Value = 42.384617
$ws.Range("J38").Value = 59095.47
$ws.Range("K38").Value = 127.153851
$ws.Range("L38").Value = 177286.41
$ws.Range("M38").Value = 219.846149
$ws.Range("N38").Value = -177980.41
# Row 44
$ws.Range("H44").Value = 5333666
$ws.Range("J44").Value = 8000372.5
$ws.Range("L44").Value = 24001117.5
$ws.Range("N44").Value = -24001913.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 135
$ws.Range("H135").Value = 35489.09
$ws.Range("J135").Value = 35489.09
$ws.Range("L135").Value = 35489.09
$ws.Range("N135").Value = -45629.09

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 108
$ws.Range("H108").Value = 32990
$ws.Range("J108").Value = 32990
$ws.Range("L108").Value = 32990
$ws.Range("N108").Value = -40670
# Row 127
$ws.Range("H127").Value = 51720.57
$ws.Range("J127").Value = 51720.57
$ws.Range("L127").Value = 51720.57
$ws.Range("N127").Value = -61640.57
# Row 132
$ws.Range("H132").Value = 11309.322
$ws.Range("I132").Value = 14015.889
$ws.Range("J132").Value = 7561.769
$ws.Range("K132").Value = 42047.667
$ws.Range("L132").Value = 22685.307
$ws.Range("M132").Value = -39517.667
$ws.Range("N132").Value = -27745.307
# Row 133
$ws.Range("H133").Value = 51633
$ws.Range("J133").Value = 51633
$ws.Range("L133").Value = 51633
$ws.Range("N133").Value = -56693
# Row 136
$ws.Range("H136").Value = 11943096
$ws.Range("I136").Value = 52356.6
$ws.Range("J136").Value = 41669944
$ws.Range("K136").Value = 157069.8
$ws.Range("L136").Value = 125009832
$ws.Range("M136").Value = -154519.8
$ws.Range("N136").Value = -125014932

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 253060.5
$ws.Range("J136").Value = 7868.385
$ws.Range("L136").Value = 23605.155
$ws.Range("N136").Value = -28705.155

